$wb = $excel.ActiveWorkbook

# --- Add the new "addEmployeeTest" sheet after the existing sheets ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "addEmployeeTest"

# --- Populate the new sheet with data ---
$newSheet.Range("A1").Value = "Username"
$newSheet.Range("B1").Value = "Password"
$newSheet.Range("C1").Value = "First Name"

$newSheet.Range("A2").Value = "Admin"
$newSheet.Range("B2").Value = "admin123"
$newSheet.Range("C2").Value = "Ken"

$newSheet.Range("A3").Value = "Admin"
$newSheet.Range("B3").Value = "admin123"
$newSheet.Range("C3").Value = "wick"

# Column widths to match bestFit columns on the other sheets
$newSheet.Columns.Item(1).ColumnWidth = 10
$newSheet.Columns.Item(2).ColumnWidth = 9.5703125
$newSheet.Columns.Item(3).ColumnWidth = 10.5703125

# Select C3 as the active cell on the new sheet (becomes the active/tabSelected sheet)
$newSheet.Range("C3").Select()

# --- Update selection on the "validCredentialTest" sheet ---
$ws2 = $wb.Worksheets.Item("validCredentialTest")
$ws2.Range("A1:B2").Select()

# Re-activate the new sheet so it is the selected tab
$newSheet.Activate()
